{"js": "// Replace the date heading and the 25 two-digit multiplication problems\n// with their updated values, matching the target commit exactly.\n//\n// Every \"old\" string is unique in the document, so a plain search would\n// normally be safe \u2014 except some \"new\" values are identical to an\n// \"old\" value used elsewhere later in the document (e.g. \"44\u00d742=\" is\n// replaced by \"36\u00d784=\", while a different, later cell's original text\n// is itself \"36\u00d784=\" and must become \"87\u00d796=\"). Replacing in a single\n// left-to-right pass would make a later search for \"36\u00d784=\" match the\n// text we *just wrote*, not the original cell.\n//\n// To avoid that collision we do this in two passes:\n//   1. Replace every old value with a unique placeholder token.\n//   2. Replace every placeholder token with its final new value.\n// Placeholder tokens never collide with any old or new value, so each\n// pass is unambiguous regardless of ordering.\n\nconst replacements = [\n  [\"2025-01-23 Thursday\", \"2025-01-24 Friday\"],\n  [\"80\u00d723=\", \"50\u00d799=\"],\n  [\"31\u00d762=\", \"86\u00d759=\"],\n  [\"45\u00d792=\", \"46\u00d711=\"],\n  [\"69\u00d782=\", \"23\u00d770=\"],\n  [\"76\u00d746=\", \"36\u00d749=\"],\n  [\"13\u00d775=\", \"23\u00d789=\"],\n  [\"71\u00d792=\", \"37\u00d749=\"],\n  [\"61\u00d741=\", \"28\u00d720=\"],\n  [\"75\u00d751=\", \"88\u00d716=\"],\n  [\"38\u00d721=\", \"40\u00d742=\"],\n  [\"44\u00d742=\", \"36\u00d784=\"],\n  [\"88\u00d764=\", \"40\u00d725=\"],\n  [\"67\u00d764=\", \"36\u00d793=\"],\n  [\"35\u00d733=\", \"40\u00d758=\"],\n  [\"42\u00d778=\", \"98\u00d787=\"],\n  [\"51\u00d766=\", \"76\u00d762=\"],\n  [\"67\u00d720=\", \"32\u00d772=\"],\n  [\"90\u00d768=\", \"95\u00d758=\"],\n  [\"57\u00d781=\", \"50\u00d783=\"],\n  [\"36\u00d784=\", \"87\u00d796=\"],\n  [\"90\u00d750=\", \"85\u00d754=\"],\n  [\"98\u00d746=\", \"83\u00d785=\"],\n  [\"53\u00d771=\", \"79\u00d753=\"],\n  [\"48\u00d728=\", \"90\u00d773=\"],\n  [\"36\u00d767=\", \"68\u00d781=\"],\n];\n\n// Pass 1: old text -> unique placeholder (document order is irrelevant\n// here because placeholders cannot collide with any old/new value).\nfor (let idx = 0; idx < replacements.length; idx++) {\n  const oldText = replacements[idx][0];\n  const placeholder = \"\\u0001PLACEHOLDER_\" + idx + \"\\u0001\";\n\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match for '\" + oldText + \"', found \" + results.items.length);\n  }\n\n  results.items[0].insertText(placeholder, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Pass 2: placeholder -> final new text.\nfor (let idx = 0; idx < replacements.length; idx++) {\n  const newText = replacements[idx][1];\n  const placeholder = \"\\u0001PLACEHOLDER_\" + idx + \"\\u0001\";\n\n  const results = context.document.body.search(placeholder, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match for placeholder \" + idx + \", found \" + results.items.length);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and the 25 two-digit multiplication problems\n# with their updated values, matching the target commit exactly.\n#\n# Every \"old\" string is unique in the document, so a plain Find/Replace\n# would normally be safe -- except some \"new\" values are identical to\n# an \"old\" value used elsewhere later in the document (e.g. \"44x42=\"\n# becomes \"36x84=\", while a different, later cell's original text is\n# itself \"36x84=\" and must become \"87x96=\"). Replacing left-to-right in\n# a single pass would make a later search for \"36x84=\" match the text\n# we *just wrote* (Find searches/wraps over the whole story), not the\n# original cell.\n#\n# To avoid that collision we do this in two passes:\n#   1. Replace every old value with a unique placeholder token built\n#      from a control character that cannot appear in any old/new\n#      value, so it can never collide.\n#   2. Replace every placeholder token with its final new value.\n\n$d = $word.ActiveDocument\n\n$mark = [char]1\n\n$oldValues = @(\n  \"2025-01-23 Thursday\",\n  \"80\u00d723=\",\n  \"31\u00d762=\",\n  \"45\u00d792=\",\n  \"69\u00d782=\",\n  \"76\u00d746=\",\n  \"13\u00d775=\",\n  \"71\u00d792=\",\n  \"61\u00d741=\",\n  \"75\u00d751=\",\n  \"38\u00d721=\",\n  \"44\u00d742=\",\n  \"88\u00d764=\",\n  \"67\u00d764=\",\n  \"35\u00d733=\",\n  \"42\u00d778=\",\n  \"51\u00d766=\",\n  \"67\u00d720=\",\n  \"90\u00d768=\",\n  \"57\u00d781=\",\n  \"36\u00d784=\",\n  \"90\u00d750=\",\n  \"98\u00d746=\",\n  \"53\u00d771=\",\n  \"48\u00d728=\",\n  \"36\u00d767=\"\n)\n\n$newValues = @(\n  \"2025-01-24 Friday\",\n  \"50\u00d799=\",\n  \"86\u00d759=\",\n  \"46\u00d711=\",\n  \"23\u00d770=\",\n  \"36\u00d749=\",\n  \"23\u00d789=\",\n  \"37\u00d749=\",\n  \"28\u00d720=\",\n  \"88\u00d716=\",\n  \"40\u00d742=\",\n  \"36\u00d784=\",\n  \"40\u00d725=\",\n  \"36\u00d793=\",\n  \"40\u00d758=\",\n  \"98\u00d787=\",\n  \"76\u00d762=\",\n  \"32\u00d772=\",\n  \"95\u00d758=\",\n  \"50\u00d783=\",\n  \"87\u00d796=\",\n  \"85\u00d754=\",\n  \"83\u00d785=\",\n  \"79\u00d753=\",\n  \"90\u00d773=\",\n  \"68\u00d781=\"\n)\n\nfunction Replace-FirstMatch($doc, $findText, $replaceText) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# Pass 1: old text -> unique placeholder (order doesn't matter here\n# because placeholders cannot collide with any old/new value).\nfor ($i = 0; $i -lt $oldValues.Length; $i++) {\n    $placeholder = $mark + \"PLACEHOLDER_\" + $i + $mark\n    Replace-FirstMatch $d $oldValues[$i] $placeholder\n}\n\n# Pass 2: placeholder -> final new text.\nfor ($i = 0; $i -lt $newValues.Length; $i++) {\n    $placeholder = $mark + \"PLACEHOLDER_\" + $i + $mark\n    Replace-FirstMatch $d $placeholder $newValues[$i]\n}\n"}
